# Apply crypto price/volume updates per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.367.23"
$ws.Range("E2").Value = "  +0.61%  "

# Row 3
$ws.Range("D3").Value = "2.488.01"
$ws.Range("E3").Value = "  +0.84%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.83%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.38%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.555"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.29%  "

# Row 9
$ws.Range("D9").Value = "2.524.06"
$ws.Range("E9").Value = "  +2.21%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0977"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.49%  "

# Row 11
$ws.Range("E11").Value = "  +0.02%  "

# Row 12
$ws.Range("E12").Value = "  -2.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.336"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.75%  "

# Row 14
$ws.Range("D14").Value = "2.941.47"
$ws.Range("E14").Value = "  +1.28%  "

# Row 15
$ws.Range("D15").Value = "58.269.00"
$ws.Range("E15").Value = "  +0.57%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.28"
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.38%  "

# Row 18
$ws.Range("D18").Value = "2.519.16"
$ws.Range("E18").Value = "  +2.41%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "325.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.72%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.58%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.37%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.51%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.408"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.61%  "

# Row 26
$ws.Range("E26").Value = "  +1.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.988"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.03%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0751"
$ws.Range("E29").Value = "  +0.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.06%  "

# Row 31
$ws.Range("E31").Value = "  +1.26%  "

# Row 32
$ws.Range("E32").Value = "  +3.14%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.13%  "

# Row 34
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.993"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.61%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.30%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.98%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.28%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.29%  "

# Row 40
$ws.Range("E40").Value = "  -0.82%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.785"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.34%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "279.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.60%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.87%  "

# Row 44
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.600"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.48%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0922"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.96%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.97%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0500"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.74%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.21%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0214"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.62%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.60%  "
